$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.343522071838379
$ws.Range("B1").Value = 1.871571898460388
$ws.Range("C1").Value = 3.326613903045654
$ws.Range("D1").Value = 3.818994045257568
$ws.Range("E1").Value = 1.054131746292114
